$d = $word.ActiveDocument
$sec = $d.Sections(1)

# --- First-page footer (footer1.xml): Pearson logo docPr id="3" ---
# Footers(2) addresses the first-page footer story in this document.
$footerFirst = $sec.Footers(2)
[void]$footerFirst.Range.InlineShapes.Item(1).Select()
$word.Selection.InlineShapes.Item(1).Name = "image2.png"

# --- Default footer (footer2.xml): Pearson logo docPr id="2" ---
# Footers(1) addresses the default (non-first-page) footer story.
$footerDefault = $sec.Footers(1)
[void]$footerDefault.Range.InlineShapes.Item(1).Select()
$word.Selection.InlineShapes.Item(1).Name = "image2.png"

# --- First-page header (header1.xml): BTEC logo docPr id="1" ---
# Headers(2) addresses the first-page header story in this document.
$headerFirst = $sec.Headers(2)
[void]$headerFirst.Range.InlineShapes.Item(1).Select()
$word.Selection.InlineShapes.Item(1).Name = "image1.jpg"
